# feat: add 2022-Q1 data
#
# 1) The former "总计" sheet (last sheet, sheetId=6) is renamed to "2022-Q1"
#    and filled with the Q1-2022 fund holding detail rows (same shape as the
#    other quarterly sheets).
# 2) A brand-new "总计" sheet is appended at the end (it naturally receives
#    sheetId=7) and is filled with the summary table, now including the new
#    2022-Q1 row at the top.

$wb = $excel.ActiveWorkbook

# A formatting/template sheet: any of the existing quarterly detail sheets
# carries the exact sheetPr / header / column styling (style index 2: bold,
# thin border, centered) that the new sheets need.
$template = $wb.Worksheets.Item(5)   # "2021-Q4"

# ---------------------------------------------------------------------
# Step 1: repurpose the old "总计" sheet as the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(6)
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Pull over the header-row / first-column formatting from the template sheet.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows: code, name, scale, stock-position, position-ratio, market-value, rank
$q1Rows = @(
    @("012148", "国投瑞银产业趋势混合型证券投资基金A", "45.53", "92.28", "4.95", "2.2537", 6),
    @("001704", "国投瑞银进宝灵活配置混合",             "33.25", "92.49", "6.47", "2.1513", 5),
    @("012149", "国投瑞银产业趋势混合型证券投资基金C", "18.42", "92.28", "4.95", "0.9118", 6),
    @("008085", "海富通先进制造股票A",                 "3.57",  "93.07", "3.43", "0.1225", 8),
    @("009025", "海富通科技创新混合A",                 "3.17",  "92.84", "3.67", "0.1163", 8),
    @("008084", "海富通先进制造股票C",                 "1.50",  "93.07", "3.43", "0.0514", 8),
    @("009024", "海富通科技创新混合C",                 "1.02",  "92.84", "3.67", "0.0374", 8),
    @("004557", "北信瑞丰鼎丰灵活配置混合",             "0.39",  "64.13", "5.15", "0.0201", 4)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = ($r - 2)
    $q1.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1.Cells.Item($r, 6).Value = "'" + $row[4]
    $q1.Cells.Item($r, 7).Value = "'" + $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet after "2022-Q1" with the updated
# summary table (adds the 2022-Q1 row on top, re-numbering column A).
# Cloning the template sheet keeps the sheetPr/outlinePr block intact.
# ---------------------------------------------------------------------
$template.Copy($null, $q1)
$total = $wb.Worksheets.Item($q1.Index + 1)
$total.Name = "总计"
$total.Cells.Clear()

$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 8,  5.66),
    @("2021-Q4", 63, 16.5),
    @("2021-Q3", 13, 6.26),
    @("2021-Q2", 25, 14.11),
    @("2021-Q1", 10, 7.86),
    @("2020-Q4", 1,  0.08)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = ($r - 2)
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

# Restore the original active sheet/tab selection (the first sheet was the
# active one before this edit).
$wb.Worksheets.Item(1).Activate()

Write-Output "done"
